# Add a second "тема2" (topic 2) row to the active sheet (Лист7), reusing
# topic 1's material/test links, matching the "new versions 17 jan 22" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing rows 3 & 4 belong to "тема1" (топик 1) but were missing the
# repeated label in column A — fill it in, matching row 2.
$ws.Range("A3").Value = "тема1"
$ws.Range("A4").Value = "тема1"

# New row 5: topic 2, reusing the same material/test reference text as
# topic 1's first row.
$ws.Range("A5").Value = "тема2"
$ws.Range("C5").Value = "материал1"
$ws.Range("D5").Value = "ссылка на материал1"
$ws.Range("E5").Value = "тест1"
$ws.Range("F5").Value = "Сыллка на тест1"

# Move the active selection to the newly-added row, as it was left after
# the edit.
$null = $ws.Range("A5").Select()
